$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 172 (pushes existing rows 172:202 down to 173:203,
# carrying their values/formatting with them).
$ws.Rows(172).Insert()

# Populate the new row 172 with the new data point (weekly Fruta/Hortalizas entry).
$ws.Cells.Item(172, 1).Value  = 11
$ws.Cells.Item(172, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(172, 3).Value  = "Bíobío"
$ws.Cells.Item(172, 4).Value  = 45244
$ws.Cells.Item(172, 5).Value  = 8
$ws.Cells.Item(172, 6).Value  = "Fruta"
$ws.Cells.Item(172, 7).Value  = 100108
$ws.Cells.Item(172, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(172, 9).Value  = 100108002
$ws.Cells.Item(172, 10).Value = "Mango"
$ws.Cells.Item(172, 11).Value = "Sin especificar"
$ws.Cells.Item(172, 12).Value = "Primera"
$ws.Cells.Item(172, 13).Value = 200
$ws.Cells.Item(172, 14).Value = 12000
$ws.Cells.Item(172, 15).Value = 13000
$ws.Cells.Item(172, 16).Value = 12500
$ws.Cells.Item(172, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(172, 18).Value = "Brasil"
$ws.Cells.Item(172, 19).Value = 3125
$ws.Cells.Item(172, 20).Value = 4
